# Fruta / hortaliza, semanal
# Insert a new weekly price-observation row at row 206 (pushing the
# existing rows 206:262 down to 207:263), then populate it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 206.
$ws.Rows.Item(206).Insert()

# Populate the new row 206 with the new observation.
$ws.Range("A206").Value = 11
$ws.Range("B206").Value = 'Vega Monumental Concepción'
$ws.Range("C206").Value = 'Bíobío'
$ws.Range("D206").Value = 44932
$ws.Range("E206").Value = 8
$ws.Range("F206").Value = 100112040
$ws.Range("G206").Value = 'Cilantro'
$ws.Range("H206").Value = 'Sin especificar'
$ws.Range("I206").Value = 'Primera'
$ws.Range("J206").Value = 270
$ws.Range("K206").Value = 18000
$ws.Range("L206").Value = 19000
$ws.Range("M206").Value = 18556
$ws.Range("N206").Value = '$/caja 36 atados'
$ws.Range("O206").Value = 'Región Metropolitana'
$ws.Range("P206").Value = 515
$ws.Range("Q206").Value = 36
$ws.Range("R206").Value = 'Hortaliza'
